# Update the "Sheet1" worksheet (the active tab) so that cell A3, which
# currently holds "Ustron, Poland", instead holds "USA". This introduces a
# new shared-string entry and also collapses the current row-wide selection
# down to the single cell A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Activate()
$ws.Range("A3").Value = "USA"
$ws.Range("A3").Select()
